# Transferencia.xlsx - add two new data-driven test scenarios
# (pagos / transferencias desde cuenta sin saldo) as a new row 9 in the
# "Datos" sheet, copying the formatting of the last existing scenario
# (row 8) and giving it the next sequential ID.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 is the last populated data row (ID 7). Duplicate it into row 9,
# keeping every column's formatting/number format identical, then bump
# the ID column to 8 - this is how a new test case was appended to the
# data sheet.
[void]$ws.Range("A8:T8").Copy($ws.Range("A9:T9"))
$ws.Range("A9").Value = 8

$excel.CutCopyMode = $false

# Leave the selection where the author left it when they saved the file.
[void]$ws.Range("A10").Select()
